$wb = $excel.ActiveWorkbook

$wsFolder  = $wb.Worksheets.Item("Folder Inventory")
$wsMeta    = $wb.Worksheets.Item("Metadata")
$wsSummary = $wb.Worksheets.Item("Summary")

# --- Folder Inventory sheet: the two most-recent folders swapped order,
#     and the newly-top folder picked up a fresh "Last Updated" timestamp ---
$wsFolder.Range("A2").Value = "Getting_started_with_Azure_AI_services"
$wsFolder.Range("B2").Value = "Getting_started_with_Azure_AI_services"
$wsFolder.Range("C2").Value = "2025-06-16 11:01:18 +0530"

$wsFolder.Range("A3").Value = "Power_Platform_Workshop:Administration_and_Governance"
$wsFolder.Range("B3").Value = "Power_Platform_Workshop:Administration_and_Governance"
$wsFolder.Range("C3").Value = "2025-06-16 10:59:35 +0530"

# --- Metadata sheet: regeneration timestamp + workflow run number ---
$wsMeta.Range("B3").Value = "2025-06-16 05:31:35 UTC"
$wsMeta.Range("B5").Value = "'11"

# --- Summary sheet: most recent update timestamp ---
$wsSummary.Range("B5").Value = "2025-06-16 11:01:18 +0530"
